$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"
$ws.Range("D2").Value = "-"
$ws.Range("C3").Value = "Desenho Técnico"
$ws.Range("C4").Value = "Desenho Técnico"
